$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A2").NumberFormat = "mm-dd-yy"
$ws.Range("A2").Value = "02/07/2023"
$ws.Range("B2").Value = "spring basics,login page"

$ws.Range("A2").Copy()
$ws.Range("A3").PasteSpecial(-4122)
$ws.Range("A3").Value = "02/08/2023"
$ws.Range("B3").Value = "c++ exception handling, working on login page"

$ws.Columns.Item(1).ColumnWidth = 9.5

$ws.Range("C6").Select()
